# Reorder the worksheet tabs: move "Summary" so it becomes the first tab,
# i.e. before "Combined Tables with Summary". "about this data" keeps its
# relative position after "Combined Tables with Summary", so the final tab
# order becomes: Summary, Combined Tables with Summary, about this data.
$wb = $excel.ActiveWorkbook

$summarySheet  = $wb.Worksheets.Item("Summary")
$combinedSheet = $wb.Worksheets.Item("Combined Tables with Summary")
$summarySheet.Move($combinedSheet)

# Re-fetch fresh references to the sheets (by name) now that the tab order
# has changed, then update the selected cell on each sheet and make
# "Summary" the active tab, matching the saved view state of the workbook.
$combinedSheet = $wb.Worksheets.Item("Combined Tables with Summary")
$combinedSheet.Activate()
$combinedSheet.Range("G4").Select() | Out-Null

$summarySheet = $wb.Worksheets.Item("Summary")
$summarySheet.Activate()
$summarySheet.Range("B4").Select() | Out-Null
